$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19
$ws.Range("A19").Value = "2NLFB6"
$ws.Range("B19").Value = "2024-12-25 19:29:21"
$ws.Range("C19").Value = "POST /setup"
$ws.Range("D19").Value = 400
$ws.Range("E19").Value = $false
$ws.Range("F19").Value = "Model Embedder untuk 'openai' harus salah satu dari ['text-embedding-3-large', 'text-embedding-3-small']."

# Row 20
$ws.Range("A20").Value = "3U3J8J"
$ws.Range("B20").Value = "2024-12-25 19:29:28"
$ws.Range("C20").Value = "POST /setup"
$ws.Range("D20").Value = 200
$ws.Range("E20").Value = $true

$f20 = @"
Proses penyiapan dokumen berhasil diselesaikan dan embeddings berhasil disimpan pada vector database.
###
llm:openai
###
model_llm:gpt-4o
###
embbeder:openai
###
model_embedder:text-embedding-3-large
###
chunk_size:900
###
chunk_overlap:100
###
total_chunks:177
"@
$ws.Range("F20").Value = $f20
